$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match column G's width to column F's width (new "STATUS" column)
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Add the new "STATUS" header in G1, matching the style/format used by the
# other header cells (copy F1's formatting onto G1)
$ws.Range("G1").Value = "STATUS"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selected cell as recorded in the saved view state
$ws.Range("M9").Select() | Out-Null
